$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 733.86884
$ws.Range("J17").Value = 610.88464
$ws.Range("L17").Value = 1832.65392
$ws.Range("N17").Value = -2168.65392

$ws.Range("H32").Value = 7500
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -4674
$ws.Range("N32").Value = -10652

$ws.Range("H44").Value = 17782.857
$ws.Range("J44").Value = 17782.857
$ws.Range("L44").Value = 17782.857
$ws.Range("N44").Value = -18706.857

$ws.Range("H75").Value = 26438.334
$ws.Range("J75").Value = 38657.5
$ws.Range("L75").Value = 38657.5
$ws.Range("N75").Value = -40529.5

$ws.Range("H78").Value = 26438.334
$ws.Range("J78").Value = 38657.5
$ws.Range("L78").Value = 115972.5
$ws.Range("N78").Value = -125332.5

$ws.Range("H97").Value = 1333.3334
$ws.Range("I97").Value = 2000
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 6000
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -5504
$ws.Range("N97").Value = -3992

$ws.Range("H137").Value = 2547.1892
$ws.Range("I137").Value = 1316.7916
$ws.Range("K137").Value = 3950.3748
$ws.Range("M137").Value = -1400.3748

$ws.Range("H138").Value = 5445.8184
$ws.Range("I138").Value = 867.4
$ws.Range("J138").Value = 7436.4346
$ws.Range("K138").Value = 2602.2
$ws.Range("L138").Value = 22309.3038
$ws.Range("M138").Value = 2537.8
$ws.Range("N138").Value = -32589.3038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3451.111
$ws.Range("I122").Value = 1765
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 5295
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -2845
$ws.Range("N122").Value = -19300

$ws.Range("H132").Value = 1950.4722
$ws.Range("I132").Value = 931.53845
$ws.Range("J132").Value = 4599.7
$ws.Range("K132").Value = 2794.61535
$ws.Range("L132").Value = 13799.1
$ws.Range("M132").Value = -264.61535
$ws.Range("N132").Value = -18859.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 6865.375
$ws.Range("I102").Value = 6865.375
$ws.Range("K102").Value = 6865.375
$ws.Range("M102").Value = -3620.375

$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws.Range("H118").Value = 28890
$ws.Range("J118").Value = 28890
$ws.Range("L118").Value = 28890
$ws.Range("N118").Value = -32204

$ws.Range("H134").Value = 1419.6981
$ws.Range("I134").Value = 957.5333000000001
$ws.Range("J134").Value = 4019.375
$ws.Range("K134").Value = 2872.5999
$ws.Range("L134").Value = 12058.125
$ws.Range("M134").Value = -337.5999000000002
$ws.Range("N134").Value = -17128.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13891026
$ws.Range("I31").Value = 826.7826
$ws.Range("J31").Value = 38465996
$ws.Range("K31").Value = 826.7826
$ws.Range("L31").Value = 38465996
$ws.Range("M31").Value = -531.7826
$ws.Range("N31").Value = -38466586

$ws.Range("H34").Value = 13891026
$ws.Range("I34").Value = 826.7826
$ws.Range("J34").Value = 38465996
$ws.Range("K34").Value = 826.7826
$ws.Range("L34").Value = 38465996
$ws.Range("M34").Value = -624.7826
$ws.Range("N34").Value = -38466400

$ws.Range("H58").Value = 2020.877
$ws.Range("I58").Value = 1771.193
$ws.Range("J58").Value = 3799.875
$ws.Range("K58").Value = 1771.193
$ws.Range("L58").Value = 3799.875
$ws.Range("M58").Value = -1568.193
$ws.Range("N58").Value = -4205.875

$ws.Range("H68").Value = 53282.125
$ws.Range("J68").Value = 53282.125
$ws.Range("L68").Value = 53282.125
$ws.Range("N68").Value = -54780.125

$ws.Range("H71").Value = 53282.125
$ws.Range("J71").Value = 53282.125
$ws.Range("L71").Value = 159846.375
$ws.Range("N71").Value = -167334.375

$ws.Range("H81").Value = 35800
$ws.Range("J81").Value = 35800
$ws.Range("L81").Value = 35800
$ws.Range("N81").Value = -37796

$ws.Range("H84").Value = 35800
$ws.Range("J84").Value = 35800
$ws.Range("L84").Value = 107400
$ws.Range("N84").Value = -117384

$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524

$ws.Range("H122").Value = 4324.4
$ws.Range("I122").Value = 3257.5
$ws.Range("J122").Value = 5924.75
$ws.Range("K122").Value = 9772.5
$ws.Range("L122").Value = 17774.25
$ws.Range("M122").Value = -7322.5
$ws.Range("N122").Value = -22674.25

$ws.Range("H123").Value = 39468
$ws.Range("J123").Value = 39468
$ws.Range("L123").Value = 39468
$ws.Range("N123").Value = -49268

$ws.Range("H132").Value = 3382.7317
$ws.Range("I132").Value = 2778.4285
$ws.Range("J132").Value = 4684.3076
$ws.Range("K132").Value = 8335.2855
$ws.Range("L132").Value = 14052.9228
$ws.Range("M132").Value = -5805.2855
$ws.Range("N132").Value = -19112.9228

$ws.Range("H134").Value = 4786.7334
$ws.Range("I134").Value = 5269.5
$ws.Range("J134").Value = 3459.125
$ws.Range("K134").Value = 15808.5
$ws.Range("L134").Value = 10377.375
$ws.Range("M134").Value = -13273.5
$ws.Range("N134").Value = -15447.375

$ws.Range("H136").Value = 2020.877
$ws.Range("I136").Value = 1771.193
$ws.Range("J136").Value = 3799.875
$ws.Range("K136").Value = 5313.579
$ws.Range("L136").Value = 11399.625
$ws.Range("M136").Value = -2763.579
$ws.Range("N136").Value = -16499.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 7385
$ws.Range("I69").Value = 1820
$ws.Range("K69").Value = 5460
$ws.Range("M69").Value = -4649

$ws.Range("H72").Value = 7385
$ws.Range("I72").Value = 1820
$ws.Range("K72").Value = 16380
$ws.Range("M72").Value = -12324

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 11158.111
$ws.Range("I41").Value = 2840.6
$ws.Range("J41").Value = 21555
$ws.Range("K41").Value = 2840.6
$ws.Range("L41").Value = 21555
$ws.Range("M41").Value = -2485.6
$ws.Range("N41").Value = -22265

$ws.Range("H122").Value = 2098.6086
$ws.Range("I122").Value = 1536.6666
$ws.Range("J122").Value = 4121.6
$ws.Range("K122").Value = 4609.9998
$ws.Range("L122").Value = 12364.8
$ws.Range("M122").Value = -2159.9998
$ws.Range("N122").Value = -17264.8

$ws.Range("H132").Value = 3174.0952
$ws.Range("I132").Value = 1775.6428
$ws.Range("J132").Value = 5971
$ws.Range("K132").Value = 5326.928400000001
$ws.Range("L132").Value = 17913
$ws.Range("M132").Value = -2796.928400000001
$ws.Range("N132").Value = -22973

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3928.1538
$ws.Range("I7").Value = 2452.2856
$ws.Range("J7").Value = 5650
$ws.Range("K7").Value = 2452.2856
$ws.Range("L7").Value = 5650
$ws.Range("M7").Value = -2340.2856
$ws.Range("N7").Value = -5874

$ws.Range("H74").Value = 40429.25
$ws.Range("J74").Value = 44776.285
$ws.Range("L74").Value = 44776.285
$ws.Range("N74").Value = -46772.285

$ws.Range("H77").Value = 40429.25
$ws.Range("J77").Value = 44776.285
$ws.Range("L77").Value = 134328.855
$ws.Range("N77").Value = -144312.855

$ws.Range("H96").Value = 30197
$ws.Range("J96").Value = 30197
$ws.Range("L96").Value = 30197
$ws.Range("N96").Value = -35689

$ws.Range("H112").Value = 31710.525
$ws.Range("J112").Value = 31710.525
$ws.Range("L112").Value = 31710.525
$ws.Range("N112").Value = -34664.525

$ws.Range("H126").Value = 3928.1538
$ws.Range("I126").Value = 2452.2856
$ws.Range("J126").Value = 5650
$ws.Range("K126").Value = 7356.8568
$ws.Range("L126").Value = 16950
$ws.Range("M126").Value = -4886.8568
$ws.Range("N126").Value = -21890

$ws.Range("H132").Value = 3334.8958
$ws.Range("I132").Value = 1468.8485
$ws.Range("J132").Value = 7440.2
$ws.Range("K132").Value = 4406.5455
$ws.Range("L132").Value = 22320.6
$ws.Range("M132").Value = -1876.5455
$ws.Range("N132").Value = -27380.6

$ws.Range("H136").Value = 2241.525
$ws.Range("I136").Value = 1298.7
$ws.Range("J136").Value = 5070
$ws.Range("K136").Value = 3896.1
$ws.Range("L136").Value = 15210
$ws.Range("M136").Value = -1346.1
$ws.Range("N136").Value = -20310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3312.8096
$ws.Range("I122").Value = 2010.5625
$ws.Range("J122").Value = 7480
$ws.Range("K122").Value = 6031.6875
$ws.Range("L122").Value = 22440
$ws.Range("M122").Value = -3581.6875
$ws.Range("N122").Value = -27340

$ws.Range("H132").Value = 8335198.5
$ws.Range("I132").Value = 1101.3462
$ws.Range("J132").Value = 23812808
$ws.Range("K132").Value = 3304.0386
$ws.Range("L132").Value = 71438424
$ws.Range("M132").Value = -774.0385999999999
$ws.Range("N132").Value = -71443484

$ws.Range("H136").Value = 2991.3103
$ws.Range("I136").Value = 844.2
$ws.Range("J136").Value = 7762.6665
$ws.Range("K136").Value = 2532.6
$ws.Range("L136").Value = 23287.9995
$ws.Range("M136").Value = 17.39999999999964
$ws.Range("N136").Value = -28387.9995
